# Fruta / hortaliza, semanal
# Update Fecha (D), Volumen (M) and, for rows 3 & 10, the unit of
# commercialization (Q), Precio $/Kg (S) and Kg/unidad (T) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44327
$ws.Range("M2").Value = 60

# Row 3
$ws.Range("D3").Value = 44322
$ws.Range("M3").Value = 60
$ws.Range("Q3").Value = "$/caja 10 kilos empedrada"
$ws.Range("S3").Value = 11500
$ws.Range("T3").Value = 1

# Row 4
$ws.Range("D4").Value = 44302
$ws.Range("M4").Value = 80

# Row 5
$ws.Range("D5").Value = 44316
$ws.Range("M5").Value = 120

# Row 6
$ws.Range("D6").Value = 44330
$ws.Range("M6").Value = 60

# Row 7
$ws.Range("D7").Value = 44306

# Row 8
$ws.Range("D8").Value = 44313
$ws.Range("M8").Value = 120

# Row 9
$ws.Range("D9").Value = 44323
$ws.Range("M9").Value = 80

# Row 10
$ws.Range("D10").Value = 44309
$ws.Range("M10").Value = 80
$ws.Range("Q10").Value = "$/caja 14 kilos granel"
$ws.Range("S10").Value = 821
$ws.Range("T10").Value = 14
